$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3626.8538
$ws.Range("I40").Value = 4480.033
$ws.Range("J40").Value = 1300
$ws.Range("K40").Value = 4480.033
$ws.Range("L40").Value = 1300
$ws.Range("M40").Value = -4305.033
$ws.Range("N40").Value = -1650
$ws.Range("H48").Value = 3166.6667
$ws.Range("J48").Value = 4000
$ws.Range("L48").Value = 12000
$ws.Range("N48").Value = -12584
$ws.Range("H56").Value = 3166.6667
$ws.Range("J56").Value = 4000
$ws.Range("L56").Value = 12000
$ws.Range("N56").Value = -13068
$ws.Range("H64").Value = 23669.625
$ws.Range("I64").Value = 2722
$ws.Range("J64").Value = 42941.44
$ws.Range("K64").Value = 2722
$ws.Range("L64").Value = 42941.44
$ws.Range("M64").Value = -2474
$ws.Range("N64").Value = -43437.44
$ws.Range("H67").Value = 23669.625
$ws.Range("I67").Value = 2722
$ws.Range("J67").Value = 42941.44
$ws.Range("K67").Value = 2722
$ws.Range("L67").Value = 42941.44
$ws.Range("M67").Value = -1864
$ws.Range("N67").Value = -44657.44
$ws.Range("H74").Value = 3147.9395
$ws.Range("I74").Value = 3083.9167
$ws.Range("K74").Value = 3083.9167
$ws.Range("M74").Value = -2147.9167
$ws.Range("H77").Value = 3147.9395
$ws.Range("I77").Value = 3083.9167
$ws.Range("K77").Value = 15419.5835
$ws.Range("M77").Value = -10739.5835
$ws.Range("H87").Value = 23326.666
$ws.Range("J87").Value = 23326.666
$ws.Range("L87").Value = 23326.666
$ws.Range("N87").Value = -25822.666
$ws.Range("H90").Value = 23326.666
$ws.Range("J90").Value = 23326.666
$ws.Range("L90").Value = 69979.998
$ws.Range("N90").Value = -82459.998
$ws.Range("H116").Value = 16669546
$ws.Range("I116").Value = 1977.2222
$ws.Range("J116").Value = 41670900
$ws.Range("K116").Value = 1977.2222
$ws.Range("L116").Value = 41670900
$ws.Range("M116").Value = 1464.7778
$ws.Range("N116").Value = -41677784
$ws.Range("H121").Value = 958.7692
$ws.Range("J121").Value = 957.0833
$ws.Range("L121").Value = 2871.2499
$ws.Range("N121").Value = -6365.2499
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2094.205
$ws.Range("I63").Value = 1934.7858
$ws.Range("J63").Value = 2500
$ws.Range("K63").Value = 1934.7858
$ws.Range("L63").Value = 2500
$ws.Range("M63").Value = -1248.7858
$ws.Range("N63").Value = -3872
$ws.Range("H66").Value = 2094.205
$ws.Range("I66").Value = 1934.7858
$ws.Range("J66").Value = 2500
$ws.Range("K66").Value = 9673.929
$ws.Range("L66").Value = 12500
$ws.Range("M66").Value = -6241.929
$ws.Range("N66").Value = -19364
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 24397.166
$ws.Range("I82").Value = 2400
$ws.Range("J82").Value = 46394.332
$ws.Range("K82").Value = 2400
$ws.Range("L82").Value = 46394.332
$ws.Range("M82").Value = -2017
$ws.Range("N82").Value = -47160.332
$ws.Range("H85").Value = 24397.166
$ws.Range("I85").Value = 2400
$ws.Range("J85").Value = 46394.332
$ws.Range("K85").Value = 2400
$ws.Range("L85").Value = 46394.332
$ws.Range("M85").Value = -1074
$ws.Range("N85").Value = -49046.332
$ws.Range("H86").Value = 4376.3784
$ws.Range("I86").Value = 4075.5938
$ws.Range("J86").Value = 6301.4
$ws.Range("K86").Value = 4075.5938
$ws.Range("L86").Value = 6301.4
$ws.Range("M86").Value = -2952.5938
$ws.Range("N86").Value = -8547.4
$ws.Range("H89").Value = 4376.3784
$ws.Range("I89").Value = 4075.5938
$ws.Range("J89").Value = 6301.4
$ws.Range("K89").Value = 20377.969
$ws.Range("L89").Value = 31507
$ws.Range("M89").Value = -14761.969
$ws.Range("N89").Value = -42739
$ws.Range("H134").Value = 1352.3871
$ws.Range("I134").Value = 1231.6923
$ws.Range("J134").Value = 1980
$ws.Range("K134").Value = 3695.0769
$ws.Range("L134").Value = 5940
$ws.Range("M134").Value = -1160.0769
$ws.Range("N134").Value = -11010
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 745.4074
$ws.Range("I22").Value = 800.6667
$ws.Range("J22").Value = 634.8889
$ws.Range("K22").Value = 800.6667
$ws.Range("L22").Value = 634.8889
$ws.Range("M22").Value = -450.6667
$ws.Range("N22").Value = -1334.8889
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 717.97296
$ws.Range("I5").Value = 289.5
$ws.Range("J5").Value = 2554.2856
$ws.Range("K5").Value = 868.5
$ws.Range("L5").Value = 7662.8568
$ws.Range("M5").Value = -756.5
$ws.Range("N5").Value = -7886.8568
$ws.Range("H34").Value = 1354.1
$ws.Range("I34").Value = 1680.3334
$ws.Range("J34").Value = 1214.2858
$ws.Range("K34").Value = 5041.0002
$ws.Range("L34").Value = 3642.8574
$ws.Range("M34").Value = -4957.0002
$ws.Range("N34").Value = -3810.8574
$ws.Range("H122").Value = 1457.4242
$ws.Range("I122").Value = 1577.1111
$ws.Range("J122").Value = 1412.5416
$ws.Range("K122").Value = 14193.9999
$ws.Range("L122").Value = 12712.8744
$ws.Range("M122").Value = -11743.9999
$ws.Range("N122").Value = -17612.8744
$ws.Range("H131").Value = 902.25
$ws.Range("I131").Value = 545.7143
$ws.Range("J131").Value = 929.086
$ws.Range("K131").Value = 1637.1429
$ws.Range("L131").Value = 2787.258
$ws.Range("M131").Value = 3402.8571
$ws.Range("N131").Value = -12867.258
$ws.Range("H135").Value = 717.97296
$ws.Range("I135").Value = 289.5
$ws.Range("J135").Value = 2554.2856
$ws.Range("K135").Value = 2605.5
$ws.Range("L135").Value = 22988.5704
$ws.Range("M135").Value = -70.5
$ws.Range("N135").Value = -28058.5704
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2181.6843
$ws.Range("I40").Value = 1766.3077
$ws.Range("K40").Value = 1766.3077
$ws.Range("M40").Value = -1630.3077
$ws.Range("H46").Value = 135532.06
$ws.Range("I46").Value = 813.26666
$ws.Range("J46").Value = 279873.66
$ws.Range("K46").Value = 813.26666
$ws.Range("L46").Value = 279873.66
$ws.Range("M46").Value = -625.26666
$ws.Range("N46").Value = -280249.66
$ws.Range("H61").Value = 1609.2963
$ws.Range("I61").Value = 1818
$ws.Range("J61").Value = 1113.625
$ws.Range("K61").Value = 1818
$ws.Range("L61").Value = 1113.625
$ws.Range("M61").Value = -1616
$ws.Range("N61").Value = -1517.625
$ws.Range("H113").Value = 1609.2963
$ws.Range("I113").Value = 1818
$ws.Range("J113").Value = 1113.625
$ws.Range("K113").Value = 1818
$ws.Range("L113").Value = 1113.625
$ws.Range("M113").Value = 352
$ws.Range("N113").Value = -5453.625
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 46875364
$ws.Range("I113").Value = 71428910
$ws.Range("J113").Value = 12500396
$ws.Range("K113").Value = 214286730
$ws.Range("L113").Value = 37501188
$ws.Range("M113").Value = -214284560
$ws.Range("N113").Value = -37505528
